$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data value in F2
$ws.Range("F2").Value = 4.1

# Update selection to F8
$ws.Range("F8").Select()
